$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Merge in the missing sales-tax data: the "year" column (F) was left blank
# for rows 156-205 even though the neighboring row (155) already carried
# 2009. Fill the gap so every record in this block is tagged 2009.
for ($r = 156; $r -le 205; $r++) {
    $ws.Cells.Item($r, 6).Value = 2009
}

# Reflect the resulting scroll/selection state: the user ended up further
# down the sheet with F207 (just past the last populated row) selected.
$ws.Range("A175").Select()
$ws.Range("F207").Select()
